$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '30.623.12'
$ws.Range("E2").Value = '  +0.47%  '

$ws.Range("D3").Value = '2.095.07'
$ws.Range("E3").Value = '  +4.72%  '

$ws.Range("D4").Value = "'1.004"
$ws.Range("E4").Value = '  +0.44%  '

$ws.Range("D5").Value = "'329.56"
$ws.Range("E5").Value = '  +1.70%  '

$ws.Range("D6").Value = "'1.002"
$ws.Range("E6").Value = '  +0.31%  '

$ws.Range("D7").Value = "'0.5258"
$ws.Range("E7").Value = '  +2.84%  '

$ws.Range("D8").Value = "'0.4304"
$ws.Range("E8").Value = '  +4.03%  '

$ws.Range("D9").Value = "'0.08857"
$ws.Range("E9").Value = '  +1.57%  '

$ws.Range("D10").Value = "'46.60"
$ws.Range("E10").Value = '  +8.67%  '

$ws.Range("D11").Value = "'1.163"
$ws.Range("E11").Value = '  +2.38%  '

$ws.Range("D12").Value = "'24.44"
$ws.Range("E12").Value = '  -1.45%  '

$ws.Range("D13").Value = '2.102.20'
$ws.Range("E13").Value = '  +5.29%  '

$ws.Range("D14").Value = "'6.678"
$ws.Range("E14").Value = '  +2.09%  '

$ws.Range("D15").Value = "'7.785"
$ws.Range("E15").Value = '  +4.62%  '

$ws.Range("D16").Value = "'96.48"
$ws.Range("E16").Value = '  +2.50%  '

$ws.Range("D17").Value = "'1.004"
$ws.Range("E17").Value = '  +0.38%  '

$ws.Range("D18").Value = "'0.00001123"
$ws.Range("E18").Value = '  +0.59%  '

$ws.Range("D19").Value = "'0.06650"
$ws.Range("E19").Value = '  +2.21%  '

$ws.Range("D20").Value = "'18.85"
$ws.Range("E20").Value = '  -0.57%  '

$ws.Range("D21").Value = "'1.002"
$ws.Range("E21").Value = '  +0.22%  '

$ws.Range("D22").Value = "'6.268"
$ws.Range("E22").Value = '  +1.64%  '

$ws.Range("D23").Value = '30.687.24'
$ws.Range("E23").Value = '  +0.50%  '

$ws.Range("D24").Value = "'12.31"
$ws.Range("E24").Value = '  +4.05%  '

$ws.Range("D25").Value = '2.348.04'
$ws.Range("E25").Value = '  +5.33%  '

$ws.Range("D26").Value = "'2.286"
$ws.Range("E26").Value = '  +3.38%  '

$ws.Range("D27").Value = "'22.41"
$ws.Range("E27").Value = '  -0.13%  '

$ws.Range("D28").Value = "'2.542"
$ws.Range("E28").Value = '  +4.78%  '

$ws.Range("D29").Value = "'162.35"
$ws.Range("E29").Value = '  -0.46%  '

$ws.Range("D30").Value = "'132.62"
$ws.Range("E30").Value = '  +0.79%  '

$ws.Range("D31").Value = "'1.204"
$ws.Range("E31").Value = '  +5.35%  '

$ws.Range("D32").Value = "'0.1074"
$ws.Range("E32").Value = '  +1.93%  '

$ws.Range("D33").Value = "'6.118"
$ws.Range("E33").Value = '  +0.46%  '

$ws.Range("D34").Value = "'1.539"
$ws.Range("E34").Value = '  +15.33%  '

$ws.Range("D35").Value = "'3.834"
$ws.Range("E35").Value = '  -0.09%  '

$ws.Range("D36").Value = "'0.02582"
$ws.Range("E36").Value = '  +2.67%  '

$ws.Range("D37").Value = "'9.734"
$ws.Range("E37").Value = '  +7.81%  '

$ws.Range("D38").Value = "'5.497"
$ws.Range("E38").Value = '  +2.33%  '

$ws.Range("D39").Value = "'0.06696"
$ws.Range("E39").Value = '  +1.34%  '

$ws.Range("B40").Value = 'Algorand'
$ws.Range("C40").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D40").Value = "'0.2260"
$ws.Range("E40").Value = '  +2.62%  '

$ws.Range("B41").Value = 'Aptos'
$ws.Range("C41").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D41").Value = "'12.50"
$ws.Range("E41").Value = '  +2.62%  '

$ws.Range("D42").Value = "'0.6765"
$ws.Range("E42").Value = '  +2.06%  '

$ws.Range("D43").Value = "'1.248"
$ws.Range("E43").Value = '  +1.19%  '

$ws.Range("E44").Value = '  +0.20%  '

$ws.Range("B45").Value = 'Decentraland'
$ws.Range("C45").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D45").Value = "'0.6360"
$ws.Range("E45").Value = '  +3.18%  '

$ws.Range("B46").Value = 'EnergySwap'
$ws.Range("C46").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D46").Value = "'13.94"
$ws.Range("E46").Value = '  +1.23%  '

$ws.Range("D47").Value = "'2.205"
$ws.Range("E47").Value = '  -0.20%  '

$ws.Range("D48").Value = "'3.628"
$ws.Range("E48").Value = '  -0.93%  '

$ws.Range("D49").Value = "'1.243"
$ws.Range("E49").Value = '  -1.30%  '

$ws.Range("B50").Value = 'WEMIXTOKEN'
$ws.Range("C50").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D50").Value = "'1.195"
$ws.Range("E50").Value = '  +7.28%  '

$ws.Range("B51").Value = 'Aave'
$ws.Range("C51").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D51").Value = "'82.43"
$ws.Range("E51").Value = '  +2.52%  '
